# Fruta / hortaliza, semanal
# Insert a new weekly record as row 29 ("Primera" quality, 2022-01-27),
# pushing the existing rows 29-41 down to 30-42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank row at position 29 (shifts rows 29:41 -> 30:42,
# inherits the row-28/30 number formatting, e.g. the date style on column D).
$ws.Rows("29:29").Insert()

# Populate the freshly inserted row 29 with the new weekly observation.
$ws.Cells.Item(29, 1).Value  = 10
$ws.Cells.Item(29, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(29, 3).Value  = "La Araucanía"
$ws.Cells.Item(29, 4).Value  = 44588
$ws.Cells.Item(29, 5).Value  = 9
$ws.Cells.Item(29, 6).Value  = "Fruta"
$ws.Cells.Item(29, 7).Value  = 100107
$ws.Cells.Item(29, 8).Value  = "Otros"
$ws.Cells.Item(29, 9).Value  = 100107011
$ws.Cells.Item(29, 10).Value = "Tuna"
$ws.Cells.Item(29, 11).Value = "Sin especificar"
$ws.Cells.Item(29, 12).Value = "Primera"
$ws.Cells.Item(29, 13).Value = 50
$ws.Cells.Item(29, 14).Value = 25000
$ws.Cells.Item(29, 15).Value = 25000
$ws.Cells.Item(29, 16).Value = 25000
$ws.Cells.Item(29, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(29, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(29, 19).Value = 1562
$ws.Cells.Item(29, 20).Value = 16
